$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new tweet rows (45-53) scraped for #kuliahonline
$ws.Cells.Item(45,1).Value2 = 1337702672737059072
$ws.Cells.Item(45,2).Value2 = "Selengkapnya silakan kunjungi laman berikut:⁣`nhttps://t.co/WRfnwhRBQt⁣`n⁣`n#persma⁣`n#kavling10⁣`n#kuliahonline"
$ws.Cells.Item(45,3).Value2 = "lpmkavling10"
$ws.Cells.Item(45,4).Value2 = "Sat Dec 12 10:15:45 +0000 2020"

$ws.Cells.Item(46,1).Value2 = 1337577916314447872
$ws.Cells.Item(46,2).Value2 = "📚Promo Buku Perguruan Tinggi Erlangga🎓 `n#PromoBuku #bukuibuku #BukuKuliah #JualBuku #BeliBuku #TokoBuku #BukuMurah #DiskonBuku #BukuDiskon #KuliahOnline #Mahasiswa  #belajardirumah #ingatpesanibu #pakaimasker #CuciTanganPakaiSabun #JagaJarak `n⏩ https://t.co/eIkeMmFT2n ✔ https://t.co/9CrZWFEhQm"
$ws.Cells.Item(46,3).Value2 = "PotekantropusX"
$ws.Cells.Item(46,4).Value2 = "Sat Dec 12 02:00:00 +0000 2020"

$ws.Cells.Item(47,1).Value2 = 1337295881935020032
$ws.Cells.Item(47,2).Value2 = "Bolehkah saya menyerah sampe ke titik ini🥺`n#Kuliahonline"
$ws.Cells.Item(47,3).Value2 = "LisaNurfitri3"
$ws.Cells.Item(47,4).Value2 = "Fri Dec 11 07:19:18 +0000 2020"

$ws.Cells.Item(48,1).Value2 = 1337215527815705088
$ws.Cells.Item(48,2).Value2 = "📚Promo Buku Perguruan Tinggi Erlangga🎓 `n#PromoBuku #bukuibuku #BukuKuliah #JualBuku #BeliBuku #TokoBuku #BukuMurah #DiskonBuku #BukuDiskon #BacaBuku #KuliahOnline #BelajarDiRumah #ingatpesanibu  `n⏩ https://t.co/eIkeMmXuqX ✔ https://t.co/ONlyEasEcE"
$ws.Cells.Item(48,3).Value2 = "PotekantropusX"
$ws.Cells.Item(48,4).Value2 = "Fri Dec 11 02:00:00 +0000 2020"

$ws.Cells.Item(49,1).Value2 = 1337046532605021952
$ws.Cells.Item(49,2).Value2 = "Universitas Islam Bandung (Unisba) memperpanjang perkuliahan dalam jaringan (daring) hingga ujian tengah semester genap 2021. Gimana nih tanggapan sobat kampus?Klik link di bawah ini!!`nhttps://t.co/7T6v9I0mt6`n#unisba #universitasislambandung #kuliahonline #daring https://t.co/OCiGcYkJ1e"
$ws.Cells.Item(49,3).Value2 = "suaramahasiswa"
$ws.Cells.Item(49,4).Value2 = "Thu Dec 10 14:48:29 +0000 2020"

$ws.Cells.Item(50,1).Value2 = 1336967210498019072
$ws.Cells.Item(50,2).Value2 = "kuliah gue online lagi, barusan dapet edaran bakal online 1 semester lagi. ga bisa bayangin kuliah online satu semester dengan tugas yang bejibun. Tapi ya mau gimana lagi demi kepentingan bersama.`nSemoga semuanya cepet bener ye`nbiar gue bisa hunting cogan juga 😃👍`n#KuliahOnline"
$ws.Cells.Item(50,3).Value2 = "buciinmarklee"
$ws.Cells.Item(50,4).Value2 = "Thu Dec 10 09:33:17 +0000 2020"

$ws.Cells.Item(51,1).Value2 = 1336945957015932928
$ws.Cells.Item(51,2).Value2 = "Bisniscerita kali ini mau bagiin sesuatu .... semoga bermanfaat. #Trending #MarketingTwitter #kultwit #KuliahOnline #tips https://t.co/YGWAwWdvVz"
$ws.Cells.Item(51,3).Value2 = "aziszaenalakbar"
$ws.Cells.Item(51,4).Value2 = "Thu Dec 10 08:08:50 +0000 2020"

$ws.Cells.Item(52,1).Value2 = 1336892842946354944
$ws.Cells.Item(52,2).Value2 = "capek online school ? atau pusing krena gapaham paham matematika ? yukk jokiin ke akuu aja murmer bgt 😉😻😻 Bisa konsul jugaa materi mtk, free sampe pahaamm `n`nrules nya baca di thread bawah yaa 🤗🤗`n#matematika #jokitugas #onlineschool #kuliahonline https://t.co/ulta6KVqPd"
$ws.Cells.Item(52,3).Value2 = "brawnxugar"
$ws.Cells.Item(52,4).Value2 = "Thu Dec 10 04:37:46 +0000 2020"

$ws.Cells.Item(53,1).Value2 = 1336865606121914112
$ws.Cells.Item(53,2).Value2 = "jasa joki bagi kalian penting ga sih sebenernya ??`n`n#jokitugas #murah #Tugas #KuliahOnline #skripsi #bebanhidup"
$ws.Cells.Item(53,3).Value2 = "JasaJoki_"
$ws.Cells.Item(53,4).Value2 = "Thu Dec 10 02:49:32 +0000 2020"

# Match the author's final scroll position / selection
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H55").Select()

"Added rows 45-53"
